$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1018.2222
$ws.Range("J17").Value = 1018.2222
$ws.Range("L17").Value = 3054.6666
$ws.Range("N17").Value = -3390.6666
$ws.Range("H33").Value = 139.73914
$ws.Range("I33").Value = 106.888885
$ws.Range("K33").Value = 106.888885
$ws.Range("M33").Value = 122.111115
$ws.Range("H98").Value = 2249.7778
$ws.Range("I98").Value = 1801.1428
$ws.Range("K98").Value = 1801.1428
$ws.Range("M98").Value = -303.1428000000001
$ws.Range("H107").Value = 2090.1052
$ws.Range("I107").Value = 1586.2142
$ws.Range("K107").Value = 1586.2142
$ws.Range("M107").Value = 333.7858000000001
$ws.Range("H111").Value = 5345.8
$ws.Range("I111").Value = 5345.8
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 16037.4
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -12970.4
$ws.Range("H116").Value = 3211.1667
$ws.Range("I116").Value = 2587.0833
$ws.Range("J116").Value = 4459.3335
$ws.Range("K116").Value = 2587.0833
$ws.Range("L116").Value = 4459.3335
$ws.Range("M116").Value = 854.9167000000002
$ws.Range("N116").Value = -11343.3335
$ws.Range("H122").Value = 2249.7778
$ws.Range("I122").Value = 1801.1428
$ws.Range("K122").Value = 5403.428400000001
$ws.Range("M122").Value = -2953.428400000001
$ws.Range("H132").Value = 7413532
$ws.Range("I132").Value = 11907527
$ws.Range("J132").Value = 11657.294
$ws.Range("K132").Value = 35722581
$ws.Range("L132").Value = 34971.882
$ws.Range("M132").Value = -35720051
$ws.Range("N132").Value = -40031.882
$ws.Range("H137").Value = 1783.1666
$ws.Range("I137").Value = 2999.5
$ws.Range("J137").Value = 1175
$ws.Range("K137").Value = 8998.5
$ws.Range("L137").Value = 3525
$ws.Range("M137").Value = -6448.5
$ws.Range("N137").Value = -8625
$ws.Range("H138").Value = 778643.6
$ws.Range("J138").Value = 1063061.1
$ws.Range("L138").Value = 3189183.3
$ws.Range("N138").Value = -3199463.3

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H5").Value = 300
$ws.Range("I5").Value = 300
$ws.Range("K5").Value = 300
$ws.Range("M5").Value = -188
$ws.Range("H32").Value = 4118.4736
$ws.Range("I32").Value = 4302.4243
$ws.Range("J32").Value = 2904.4
$ws.Range("K32").Value = 4302.4243
$ws.Range("L32").Value = 2904.4
$ws.Range("M32").Value = -4015.4243
$ws.Range("N32").Value = -3478.4
$ws.Range("H45").Value = 1388.4445
$ws.Range("I45").Value = 1370.1177
$ws.Range("K45").Value = 1370.1177
$ws.Range("M45").Value = -993.1177
$ws.Range("H54").Value = 14000
$ws.Range("J54").Value = 14000
$ws.Range("L54").Value = 14000
$ws.Range("N54").Value = -15538
$ws.Range("H102").Value = 13898401
$ws.Range("I102").Value = 16677681
$ws.Range("K102").Value = 16677681
$ws.Range("M102").Value = -16676059
$ws.Range("H110").Value = 1899.1333
$ws.Range("I110").Value = 1574.9231
$ws.Range("K110").Value = 1574.9231
$ws.Range("M110").Value = 470.0769
$ws.Range("H122").Value = 1199
$ws.Range("I122").Value = 1199
$ws.Range("K122").Value = 3597
$ws.Range("M122").Value = -1147
$ws.Range("H132").Value = 3913.3
$ws.Range("I132").Value = 3522.6667
$ws.Range("J132").Value = 4499.25
$ws.Range("K132").Value = 10568.0001
$ws.Range("L132").Value = 13497.75
$ws.Range("M132").Value = -8038.000100000001
$ws.Range("N132").Value = -18557.75

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("K4").Value = 300
$ws.Range("M4").Value = -185
$ws.Range("H81").Value = 8254.571
$ws.Range("J81").Value = 8254.571
$ws.Range("L81").Value = 8254.571
$ws.Range("N81").Value = -10376.571
$ws.Range("H84").Value = 8254.571
$ws.Range("J84").Value = 8254.571
$ws.Range("L84").Value = 24763.713
$ws.Range("N84").Value = -35371.713
$ws.Range("H105").Value = 166668380
$ws.Range("I105").Value = 250001730
$ws.Range("K105").Value = 250001730
$ws.Range("M105").Value = -249999983
$ws.Range("H107").Value = 1410.75
$ws.Range("I107").Value = 973.86664
$ws.Range("K107").Value = 973.86664
$ws.Range("M107").Value = 946.13336

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 367.0909
$ws.Range("I22").Value = 353.8
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 353.8
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -3.800000000000011
$ws.Range("N22").Value = -1200
$ws.Range("H31").Value = 1579.238
$ws.Range("I31").Value = 1429.6842
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 1429.6842
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -1134.6842
$ws.Range("N31").Value = -3590
$ws.Range("H34").Value = 1579.238
$ws.Range("I34").Value = 1429.6842
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 1429.6842
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -1227.6842
$ws.Range("N34").Value = -3404
$ws.Range("H107").Value = 676.5
$ws.Range("I107").Value = 436
$ws.Range("J107").Value = 796.75
$ws.Range("K107").Value = 436
$ws.Range("L107").Value = 796.75
$ws.Range("M107").Value = 1484
$ws.Range("N107").Value = -4636.75
$ws.Range("H134").Value = 27780074
$ws.Range("I134").Value = 41668924
$ws.Range("J134").Value = 2375
$ws.Range("K134").Value = 125006772
$ws.Range("L134").Value = 7125
$ws.Range("M134").Value = -125004237
$ws.Range("N134").Value = -12195

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 296.16666
$ws.Range("I7").Value = 268.75
$ws.Range("J7").Value = 351
$ws.Range("K7").Value = 806.25
$ws.Range("L7").Value = 1053
$ws.Range("M7").Value = -694.25
$ws.Range("N7").Value = -1277

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 169.66667
$ws.Range("I2").Value = 132.5
$ws.Range("K2").Value = 132.5
$ws.Range("M2").Value = -19.5
$ws.Range("H113").Value = 1988.2222
$ws.Range("J113").Value = 2496.6
$ws.Range("L113").Value = 2496.6
$ws.Range("N113").Value = -6836.6
$ws.Range("H122").Value = 1198.8889
$ws.Range("I122").Value = 1065
$ws.Range("J122").Value = 1466.6666
$ws.Range("K122").Value = 3195
$ws.Range("L122").Value = 4399.9998
$ws.Range("M122").Value = -745
$ws.Range("N122").Value = -9299.9998
$ws.Range("H132").Value = 2530.25
$ws.Range("I132").Value = 2554.0908
$ws.Range("J132").Value = 2477.8
$ws.Range("K132").Value = 7662.2724
$ws.Range("L132").Value = 7433.400000000001
$ws.Range("M132").Value = -5132.2724
$ws.Range("N132").Value = -12493.4

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2921.4
$ws.Range("I40").Value = 2667.3333
$ws.Range("J40").Value = 3302.5
$ws.Range("K40").Value = 2667.3333
$ws.Range("L40").Value = 3302.5
$ws.Range("M40").Value = -2531.3333
$ws.Range("N40").Value = -3574.5
$ws.Range("H46").Value = 1249.9333
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 1474.9
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 1474.9
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -1850.9
$ws.Range("H68").Value = 1799.5
$ws.Range("H71").Value = 1799.5
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = 0
$ws.Range("H132").Value = 79972.46000000001
$ws.Range("I132").Value = 3257.8
$ws.Range("J132").Value = 127919.125
$ws.Range("K132").Value = 9773.400000000001
$ws.Range("L132").Value = 383757.375
$ws.Range("M132").Value = -7243.400000000001
$ws.Range("N132").Value = -388817.375

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20496
$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21716
$ws.Range("H107").Value = 314.8846
$ws.Range("I107").Value = 327.83334
$ws.Range("J107").Value = 285.75
$ws.Range("K107").Value = 983.5000200000001
$ws.Range("L107").Value = 857.25
$ws.Range("M107").Value = 936.4999799999999
$ws.Range("N107").Value = -4697.25
$ws.Range("H113").Value = 815
$ws.Range("I113").Value = 600.4
$ws.Range("K113").Value = 1801.2
$ws.Range("M113").Value = 368.8000000000002
$ws.Range("H122").Value = 7650543.5
$ws.Range("I122").Value = 8670392
$ws.Range("J122").Value = 1677.5
$ws.Range("K122").Value = 26011176
$ws.Range("L122").Value = 5032.5
$ws.Range("M122").Value = -26008726
$ws.Range("N122").Value = -9932.5
$ws.Range("H132").Value = 2591.2974
$ws.Range("I132").Value = 2403.5386
$ws.Range("J132").Value = 3035.0908
$ws.Range("K132").Value = 7210.6158
$ws.Range("L132").Value = 9105.2724
$ws.Range("M132").Value = -4680.6158
$ws.Range("N132").Value = -14165.2724
$ws.Range("H136").Value = 503.6
$ws.Range("I136").Value = 258
$ws.Range("K136").Value = 774
$ws.Range("M136").Value = 1776
